$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")

# Row 2 (Leve Item ID 5489)
$ws.Range("H2").Value = 2495.3
$ws.Range("I2").Value = 499
$ws.Range("J2").Value = 2717.111
$ws.Range("K2").Value = 499
$ws.Range("L2").Value = 2717.111
$ws.Range("M2").Value = -386
$ws.Range("N2").Value = -2943.111

# Row 9 (Leve Item ID 5487)
$ws.Range("H9").Value = 453.10526
$ws.Range("I9").Value = 314.66666
$ws.Range("J9").Value = 972.25
$ws.Range("K9").Value = 314.66666
$ws.Range("L9").Value = 972.25
$ws.Range("M9").Value = -145.66666
$ws.Range("N9").Value = -1310.25

# Row 33 (Leve Item ID 5512)
$ws.Range("H33").Value = 1585.6154
$ws.Range("I33").Value = 357.4
$ws.Range("K33").Value = 357.4
$ws.Range("M33").Value = -128.4

# Row 40 (Leve Item ID 5505)
$ws.Range("H40").Value = 8176
$ws.Range("J40").Value = 8596.299999999999
$ws.Range("L40").Value = 8596.299999999999
$ws.Range("N40").Value = -8946.299999999999

# Row 48 (Leve Item ID 4587)
$ws.Range("H48").Value = 1000
$ws.Range("I48").Value = 1000
$ws.Range("J48").Value = 0
$ws.Range("K48").Value = 3000
$ws.Range("L48").Value = 0
$ws.Range("M48").Value = -2708
$ws.Range("N48").Value = $null

# Row 53 (Leve Item ID 5479)
$ws.Range("H53").Value = 308.95
$ws.Range("I53").Value = 254.18182
$ws.Range("J53").Value = 375.8889
$ws.Range("K53").Value = 254.18182
$ws.Range("L53").Value = 375.8889
$ws.Range("M53").Value = 382.81818
$ws.Range("N53").Value = -1649.8889

# Row 56 (Leve Item ID 4587)
$ws.Range("H56").Value = 1000
$ws.Range("I56").Value = 1000
$ws.Range("J56").Value = 0
$ws.Range("K56").Value = 3000
$ws.Range("L56").Value = 0
$ws.Range("M56").Value = -2466
$ws.Range("N56").Value = $null

# Row 58 (Leve Item ID 4606)
$ws.Range("H58").Value = 346.2857
$ws.Range("I58").Value = 346.2857
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 1038.8571
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = -888.8571000000002
$ws.Range("N58").Value = $null

# Row 70 (Leve Item ID 12604)
$ws.Range("H70").Value = 5676.1665
$ws.Range("I70").Value = 3000.3333
$ws.Range("J70").Value = 6568.1113
$ws.Range("K70").Value = 9000.999899999999
$ws.Range("L70").Value = 19704.3339
$ws.Range("M70").Value = -8730.999899999999
$ws.Range("N70").Value = -20244.3339

# Row 73 (Leve Item ID 12604)
$ws.Range("H73").Value = 5676.1665
$ws.Range("I73").Value = 3000.3333
$ws.Range("J73").Value = 6568.1113
$ws.Range("K73").Value = 9000.999899999999
$ws.Range("L73").Value = 19704.3339
$ws.Range("M73").Value = -8064.999899999999
$ws.Range("N73").Value = -21576.3339

# Row 116 (Leve Item ID 27778)
$ws.Range("H116").Value = 22239.416
$ws.Range("I116").Value = 61629.668
$ws.Range("J116").Value = 9109.333000000001
$ws.Range("K116").Value = 61629.668
$ws.Range("L116").Value = 9109.333000000001
$ws.Range("M116").Value = -58187.668
$ws.Range("N116").Value = -15993.333

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")

# Row 5 (Leve Item ID 5091)
$ws.Range("H5").Value = 51.333332
$ws.Range("I5").Value = 51.333332
$ws.Range("J5").Value = 51.333332
$ws.Range("K5").Value = 51.333332
$ws.Range("L5").Value = 51.333332
$ws.Range("M5").Value = 60.666668
$ws.Range("N5").Value = -275.333332

# Row 45 (Leve Item ID 27714)
$ws.Range("H45").Value = 2761.8667
$ws.Range("I45").Value = 2202.5
$ws.Range("K45").Value = 2202.5
$ws.Range("M45").Value = -1825.5

# Row 122 (Leve Item ID 36168)
$ws.Range("H122").Value = 5970.9
$ws.Range("I122").Value = 5526.375
$ws.Range("K122").Value = 16579.125
$ws.Range("M122").Value = -14129.125

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")

# Row 4 (Leve Item ID 5091)
$ws.Range("H4").Value = 51.333332
$ws.Range("I4").Value = 51.333332
$ws.Range("J4").Value = 51.333332
$ws.Range("K4").Value = 51.333332
$ws.Range("L4").Value = 51.333332
$ws.Range("M4").Value = 63.666668
$ws.Range("N4").Value = -281.333332

# Row 22 (Leve Item ID 5092)
$ws.Range("H22").Value = 984.0952
$ws.Range("I22").Value = 857.73334
$ws.Range("K22").Value = 857.73334
$ws.Range("M22").Value = -684.73334

# Row 26 (Leve Item ID 19535)
$ws.Range("H26").Value = 50000
$ws.Range("I26").Value = 50000
$ws.Range("K26").Value = 50000
$ws.Range("M26").Value = -49708

# Row 94 (Leve Item ID 19939)
$ws.Range("H94").Value = 4999.5
$ws.Range("I94").Value = 5000
$ws.Range("J94").Value = 4999
$ws.Range("K94").Value = 5000
$ws.Range("L94").Value = 4999
$ws.Range("M94").Value = -4549
$ws.Range("N94").Value = -5901

# Row 134 (Leve Item ID 43998)
$ws.Range("H134").Value = 3894.7896
$ws.Range("I134").Value = 3000.0588
$ws.Range("K134").Value = 9000.1764
$ws.Range("M134").Value = -6465.1764

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")

# Row 7 (Leve Item ID 5361)
$ws.Range("H7").Value = 834.8077
$ws.Range("I7").Value = 822.17645
$ws.Range("K7").Value = 822.17645
$ws.Range("M7").Value = -709.17645

# Row 58 (Leve Item ID 44021)
$ws.Range("H58").Value = 2571.1738
$ws.Range("I58").Value = 2667.4546
$ws.Range("J58").Value = 2482.9167
$ws.Range("K58").Value = 2667.4546
$ws.Range("L58").Value = 2482.9167
$ws.Range("M58").Value = -2464.4546
$ws.Range("N58").Value = -2888.9167

# Row 62 (Leve Item ID 12580)
$ws.Range("H62").Value = 10875149
$ws.Range("I62").Value = 4327
$ws.Range("J62").Value = 16672921
$ws.Range("K62").Value = 4327
$ws.Range("L62").Value = 16672921
$ws.Range("M62").Value = -3703
$ws.Range("N62").Value = -16674169

# Row 65 (Leve Item ID 12580)
$ws.Range("H65").Value = 10875149
$ws.Range("I65").Value = 4327
$ws.Range("J65").Value = 16672921
$ws.Range("K65").Value = 21635
$ws.Range("L65").Value = 83364605
$ws.Range("M65").Value = -18515
$ws.Range("N65").Value = -83370845

# Row 99 (Leve Item ID 36198)
$ws.Range("H99").Value = 5127.2354
$ws.Range("I99").Value = 5255.364
$ws.Range("K99").Value = 5255.364
$ws.Range("M99").Value = -3757.364

# Row 105 (Leve Item ID 19928)
$ws.Range("H105").Value = 2380.9167
$ws.Range("I105").Value = 2415.6365
$ws.Range("J105").Value = 1999
$ws.Range("K105").Value = 2415.6365
$ws.Range("L105").Value = 1999
$ws.Range("M105").Value = -668.6365000000001
$ws.Range("N105").Value = -5493

# Row 126 (Leve Item ID 36198)
$ws.Range("H126").Value = 5127.2354
$ws.Range("I126").Value = 5255.364
$ws.Range("K126").Value = 15766.092
$ws.Range("M126").Value = -13296.092

# Row 132 (Leve Item ID 44019)
$ws.Range("H132").Value = 6210.7915
$ws.Range("I132").Value = 3966.3157
$ws.Range("K132").Value = 11898.9471
$ws.Range("M132").Value = -9368.947100000001

# Row 136 (Leve Item ID 44021)
$ws.Range("H136").Value = 2571.1738
$ws.Range("I136").Value = 2667.4546
$ws.Range("J136").Value = 2482.9167
$ws.Range("K136").Value = 8002.3638
$ws.Range("L136").Value = 7448.750100000001
$ws.Range("M136").Value = -5452.3638
$ws.Range("N136").Value = -12548.7501

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")

# Row 7 (Leve Item ID 4728)
$ws.Range("H7").Value = 54
$ws.Range("I7").Value = 54
$ws.Range("J7").Value = 54
$ws.Range("K7").Value = 162
$ws.Range("L7").Value = 162
$ws.Range("M7").Value = -50
$ws.Range("N7").Value = -386

# Row 8 (Leve Item ID 16734)
$ws.Range("H8").Value = 612.4666999999999
$ws.Range("I8").Value = 612.4666999999999
$ws.Range("K8").Value = 1837.4001
$ws.Range("M8").Value = -1698.4001

# Row 92 (Leve Item ID 19841)
$ws.Range("H92").Value = 381.58334
$ws.Range("J92").Value = 403.47058
$ws.Range("L92").Value = 1210.41174
$ws.Range("N92").Value = -3706.41174

# Row 113 (Leve Item ID 27843)
$ws.Range("H113").Value = 632.3043
$ws.Range("I113").Value = 314.9
$ws.Range("J113").Value = 876.46155
$ws.Range("K113").Value = 944.6999999999999
$ws.Range("L113").Value = 2629.38465
$ws.Range("M113").Value = 1225.3
$ws.Range("N113").Value = -6969.38465

# Row 131 (Leve Item ID 36060)
$ws.Range("H131").Value = 3176.4285
$ws.Range("J131").Value = 4104.4546
$ws.Range("L131").Value = 12313.3638
$ws.Range("N131").Value = -22393.3638

# Row 132 (Leve Item ID 43972)
$ws.Range("H132").Value = 1487
$ws.Range("I132").Value = 1411.025
$ws.Range("K132").Value = 12699.225
$ws.Range("M132").Value = -10169.225

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")

# Row 80 (Leve Item ID 12521)
$ws.Range("H80").Value = 2540.7693
$ws.Range("I80").Value = 2335.8333
$ws.Range("K80").Value = 2335.8333
$ws.Range("M80").Value = -1337.8333

# Row 83 (Leve Item ID 12521)
$ws.Range("H83").Value = 2540.7693
$ws.Range("I83").Value = 2335.8333
$ws.Range("K83").Value = 11679.1665
$ws.Range("M83").Value = -6687.166499999999

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")

# Row 46 (Leve Item ID 5282)
$ws.Range("H46").Value = 2805.55
$ws.Range("I46").Value = 996.25
$ws.Range("J46").Value = 3257.875
$ws.Range("K46").Value = 996.25
$ws.Range("L46").Value = 3257.875
$ws.Range("M46").Value = -808.25
$ws.Range("N46").Value = -3633.875

# Row 136 (Leve Item ID 44060)
$ws.Range("H136").Value = 6796.7607
$ws.Range("I136").Value = 5469.78
$ws.Range("J136").Value = 9956.237999999999
$ws.Range("K136").Value = 16409.34
$ws.Range("L136").Value = 29868.714
$ws.Range("M136").Value = -13859.34
$ws.Range("N136").Value = -34968.714

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")

# Row 62 (Leve Item ID 12589)
$ws.Range("H62").Value = 6666.1665
$ws.Range("I62").Value = 6332.3335
$ws.Range("J62").Value = 7000
$ws.Range("K62").Value = 6332.3335
$ws.Range("L62").Value = 7000
$ws.Range("M62").Value = -5708.3335
$ws.Range("N62").Value = -8248

# Row 65 (Leve Item ID 12589)
$ws.Range("H65").Value = 6666.1665
$ws.Range("I65").Value = 6332.3335
$ws.Range("J65").Value = 7000
$ws.Range("K65").Value = 31661.6675
$ws.Range("L65").Value = 35000
$ws.Range("M65").Value = -28541.6675
$ws.Range("N65").Value = -41240

# Row 107 (Leve Item ID 27746)
$ws.Range("H107").Value = 1194.625
$ws.Range("I107").Value = 1261.0646
$ws.Range("J107").Value = 1073.4706
$ws.Range("K107").Value = 3783.1938
$ws.Range("L107").Value = 3220.4118
$ws.Range("M107").Value = -1863.1938
$ws.Range("N107").Value = -7060.4118

# Row 130 (Leve Item ID 34705)
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").Value = $null

# Row 131 (Leve Item ID 34723)
$ws.Range("H131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").Value = $null

# Row 135 (Leve Item ID 42043)
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").Value = $null

# Row 136 (Leve Item ID 44031)
$ws.Range("H136").Value = 4305.712
$ws.Range("I136").Value = 2743.366
$ws.Range("K136").Value = 8230.098
$ws.Range("M136").Value = -5680.098

Write-Host "Applied Moogle_Profits updates"